$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.607.03'
$ws.Range("E2").Value = '  -3.22%  '
$ws.Range("D3").Value = '1.850.87'
$ws.Range("E3").Value = '  -3.77%  '
$ws.Range("E4").Value = '  -0.96%  '
$ws.Range("D5").Value = '''334.30'
$ws.Range("E5").Value = '  +2.46%  '
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").Value = '''0.4647'
$ws.Range("E7").Value = '  -3.59%  '
$ws.Range("D8").Value = '''0.3922'
$ws.Range("E8").Value = '  -3.67%  '
$ws.Range("D9").Value = '''46.48'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("D10").Value = '''0.07914'
$ws.Range("E10").Value = '  -3.88%  '
$ws.Range("D11").Value = '''0.9859'
$ws.Range("E11").Value = '  -2.49%  '
$ws.Range("D12").Value = '''22.25'
$ws.Range("E12").Value = '  -5.92%  '
$ws.Range("D13").Value = '1.980.34'
$ws.Range("E13").Value = '  +4.58%  '
$ws.Range("D14").Value = '''5.852'
$ws.Range("E14").Value = '  -3.79%  '
$ws.Range("D15").Value = '''7.014'
$ws.Range("D16").Value = '''0.06866'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = '''87.87'
$ws.Range("E17").Value = '  -4.23%  '
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '''0.00001007'
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").Value = '''17.10'
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '28.645.89'
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").Value = '''5.399'
$ws.Range("E23").Value = '  -5.10%  '
$ws.Range("E24").Value = '  -5.16%  '
$ws.Range("D25").Value = '2.277.74'
$ws.Range("E25").Value = '  +5.67%  '
$ws.Range("D26").Value = '''2.128'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("D27").Value = '''153.11'
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").Value = '''19.41'
$ws.Range("E28").Value = '  -3.14%  '
$ws.Range("D29").Value = '''6.100'
$ws.Range("E29").Value = '  -5.89%  '
$ws.Range("D30").Value = '''2.017'
$ws.Range("E30").Value = '  -3.96%  '
$ws.Range("D31").Value = '''117.66'
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").Value = '''0.9817'
$ws.Range("E32").Value = '  -3.53%  '
$ws.Range("D33").Value = '''0.09413'
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("D34").Value = '''5.371'
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("D36").Value = '''1.348'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("D37").Value = '''0.06150'
$ws.Range("E37").Value = '  -3.60%  '
$ws.Range("D38").Value = '''0.02201'
$ws.Range("E38").Value = '  -4.25%  '
$ws.Range("D39").Value = '''1.161'
$ws.Range("E39").Value = '  -2.43%  '
$ws.Range("E40").Value = '  -4.14%  '
$ws.Range("D41").Value = '''7.629'
$ws.Range("E41").Value = '  -3.41%  '
$ws.Range("E42").Value = '  -6.14%  '
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").Value = '''2.376'
$ws.Range("E44").Value = '  -4.02%  '
$ws.Range("D45").Value = '''1.248'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").Value = '''11.79'
$ws.Range("E46").Value = '  -5.31%  '
$ws.Range("D47").Value = '''0.5394'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("D48").Value = '''0.07160'
$ws.Range("E48").Value = '  -4.49%  '
$ws.Range("D49").Value = '''1.910'
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("D50").Value = '''114.09'
$ws.Range("D51").Value = '''42.74'
$ws.Range("E51").Value = '  +1.90%  '
